$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J16
$data = @(
    @(9, 9),
    @(3, 5),
    @(5, 7),
    @(4, 6),
    @(9, 9),
    @(3, 7),
    @(1, 3),
    @(9, 9),
    @(4, 7),
    @(5, 5),
    @(6, 6),
    @(6, 6),
    @(9, 9),
    @(4, 4),
    @(6, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
